$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("correlation_matrix")
$ws.Range("B2").Value = 0.848847310230273
$ws.Range("A3").Value = 0.848847310230273

$ws = $wb.Worksheets.Item("equilibrium_concentrations")
$ws.Range("A2").Value = 0.0000000657079255915789
$ws.Range("B2").Value = 0.000353335365591123
$ws.Range("C2").Value = 0.000000864629911626739
$ws.Range("D2").Value = 0.00000000000449732084461406
$ws.Range("E2").Value = 0.000000230346831802918
$ws.Range("A3").Value = 0.0000839967070662358
$ws.Range("B3").Value = 0.0000857807350894835
$ws.Range("C3").Value = 0.000268335056694733
$ws.Range("D3").Value = 0.00000178420821621364
$ws.Range("E3").Value = 0.000000000180192926758746
$ws.Range("A4").Value = 0.000347583618662726
$ws.Range("B4").Value = 0.000024670213965408
$ws.Range("C4").Value = 0.000319343151640457
$ws.Range("D4").Value = 0.00000878663865680689
$ws.Range("E4").Value = 0.0000000000435452411209539
$ws.Range("A5").Value = 0.000665672088035321
$ws.Range("B5").Value = 0.0000131491284715596
$ws.Range("C5").Value = 0.00032597380835523
$ws.Range("D5").Value = 0.0000171770631735829
$ws.Range("E5").Value = 0.0000000000227373398350435
$ws.Range("A6").Value = 0.00133310373218991
$ws.Range("B6").Value = 0.00000636835440170605
$ws.Range("C6").Value = 0.000316167015005058
$ws.Range("D6").Value = 0.0000333646338094201
$ws.Range("E6").Value = 0.0000000000113536644740305
$ws.Range("A7").Value = 0.00264712583306349
$ws.Range("B7").Value = 0.00000294410784641187
$ws.Range("C7").Value = 0.000290237611670279
$ws.Range("D7").Value = 0.0000608182805040659
$ws.Range("E7").Value = 0.00000000000571775330636466
$ws.Range("A8").Value = 0.00664397774570443
$ws.Range("B8").Value = 0.000000935641668413102
$ws.Range("C8").Value = 0.000231506460134426
$ws.Range("D8").Value = 0.000121757898265341
$ws.Range("E8").Value = 0.00000000000227809500026514
$ws.Range("A9").Value = 0.0134044254224832
$ws.Range("B9").Value = 0.00000034362691489237
$ws.Range("C9").Value = 0.00017153816754401
$ws.Range("D9").Value = 0.000182018205581655
$ws.Range("E9").Value = 0.00000000000112915041169726

$ws = $wb.Worksheets.Item("absorbance_calc_abs_errors")
$ws.Range("C2").Value = 1.14899994954435
$ws.Range("D2").Value = 1.5433261305592
$ws.Range("E2").Value = 1.64171769433111
$ws.Range("F2").Value = 1.68274593717371
$ws.Range("G2").Value = 1.70012665074334
$ws.Range("H2").Value = 1.70987231378232
$ws.Range("I2").Value = 1.74640288825853
$ws.Range("J2").Value = 1.78782732674432
$ws.Range("C3").Value = 2.70199995089303
$ws.Range("D3").Value = 2.38946851601394
$ws.Range("E3").Value = 2.24835480918796
$ws.Range("F3").Value = 2.23239307233761
$ws.Range("G3").Value = 2.18466740939406
$ws.Range("H3").Value = 2.11421833415463
$ws.Range("I3").Value = 2.01979962465506
$ws.Range("J3").Value = 1.9890621527603
$ws.Range("C4").Value = -0.0000000504556501024211
$ws.Range("D4").Value = 0.000326130559196391
$ws.Range("E4").Value = -0.000282305668892313
$ws.Range("F4").Value = -0.00225406282629392
$ws.Range("G4").Value = -0.000873349256657052
$ws.Range("H4").Value = 0.00587231378232156
$ws.Range("I4").Value = -0.00359711174146571
$ws.Range("J4").Value = 0.000827326744315293
$ws.Range("C5").Value = -0.0000000491069669372735
$ws.Range("D5").Value = 0.000468516013935272
$ws.Range("E5").Value = -0.00364519081203563
$ws.Range("F5").Value = 0.00439307233761443
$ws.Range("G5").Value = -0.00133259060594426
$ws.Range("H5").Value = 0.00021833415462913
$ws.Range("I5").Value = -0.000200375344935111
$ws.Range("J5").Value = 0.0000621527603048566

$ws = $wb.Worksheets.Item("absorbance_calc_rel_errors")
$ws.Range("C2").Value = 1.14899994954435
$ws.Range("D2").Value = 1.5433261305592
$ws.Range("E2").Value = 1.64171769433111
$ws.Range("F2").Value = 1.68274593717371
$ws.Range("G2").Value = 1.70012665074334
$ws.Range("H2").Value = 1.70987231378232
$ws.Range("I2").Value = 1.74640288825853
$ws.Range("J2").Value = 1.78782732674432
$ws.Range("C3").Value = 2.70199995089303
$ws.Range("D3").Value = 2.38946851601394
$ws.Range("E3").Value = 2.24835480918796
$ws.Range("F3").Value = 2.23239307233761
$ws.Range("G3").Value = 2.18466740939406
$ws.Range("H3").Value = 2.11421833415463
$ws.Range("I3").Value = 2.01979962465506
$ws.Range("J3").Value = 1.9890621527603
$ws.Range("C4").Value = -0.000000043912663274518
$ws.Range("D4").Value = 0.000211361347502522
$ws.Range("E4").Value = -0.000171927934769984
$ws.Range("F4").Value = -0.00133772274557503
$ws.Range("G4").Value = -0.000513432837540889
$ws.Range("H4").Value = 0.00344619353422627
$ws.Range("I4").Value = -0.00205549242369469
$ws.Range("J4").Value = 0.000462969638676717
$ws.Range("C5").Value = -0.0000000181743030855934
$ws.Range("D5").Value = 0.000196113861002625
$ws.Range("E5").Value = -0.00161864600889681
$ws.Range("F5").Value = 0.00197175598636195
$ws.Range("G5").Value = -0.00060960229000195
$ws.Range("H5").Value = 0.000103280110988235
$ws.Range("I5").Value = -0.0000991957153144113
$ws.Range("J5").Value = 0.0000312482455026931

$ws = $wb.Worksheets.Item("mol_ext_coefficients_calc")
$ws.Range("C2").Value = 3.62851768723924
$ws.Range("D2").Value = 2964.36394125175
$ws.Range("E2").Value = 4768.81513735542
$ws.Range("F2").Value = 5055.18778225329
$ws.Range("G2").Value = 423108.972180095
$ws.Range("C3").Value = 22.5538789403126
$ws.Range("D3").Value = 8048.51242828701
$ws.Range("E3").Value = 6303.21573312717
$ws.Range("F3").Value = 3311.39032240926
$ws.Range("G3").Value = -639364.45312061
$ws.Range("C4").Value = 2.2557199599467
$ws.Range("D4").Value = 70.1880926029722
$ws.Range("E4").Value = 12.6654507470743
$ws.Range("F4").Value = 164.100837409176
$ws.Range("G4").Value = 109207.061021345
$ws.Range("C5").Value = 1.80524215495195
$ws.Range("D5").Value = 56.1712028941527
$ws.Range("E5").Value = 10.1361011145313
$ws.Range("F5").Value = 131.329134207315
$ws.Range("G5").Value = 87397.9011910678


$ws = $wb.Worksheets.Item("constants_evaluated")
$pairs = @(
  @("B4", "4.571025390625"),
  @("C4", "0.108717276173831"),
  @("B5", "6.46953125"),
  @("C5", "0.260200371728513")
)
foreach ($p in $pairs) {
  $addr = $p[0]
  $val = $p[1]
  $r = $ws.Range($addr)
  $r.Formula = '="' + $val + '"'
  $r.Copy($r)
  $r.PasteSpecial(-4163)
}

$ws = $wb.Worksheets.Item("adj_r_squared")
$r = $ws.Range("A2")
$r.Formula = '=TEXT(0.999944,"0.000000")'
$r.Copy($r)
$r.PasteSpecial(-4163)
